$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet originally has two data rows (row 2 and row 3) inside Table1
# (A1:G3). The edit drops down to a single shipment record: row 3 is
# removed and row 2 is overwritten with the new shipment's data
# (vehicle plate, responsible person, route, cost and both dates).

# Shrink the table to its new extent first (A1:G2), then delete the
# now out-of-range row 3.
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:G2"))
$ws.Rows("3:3").Delete()

# Update the (now) row 2 values to the new shipment's data.
$ws.Range("A2").Value = "FMG399"
$ws.Range("B2").Value = "PEDRO PEREZ"
$ws.Range("C2").Value = " CEMENTO, LADRILLOS"
$ws.Range("D2").Value = "SANTANDER / BUCARAMANGA    -    ANTIOQUIA / MEDELLIN"
$ws.Range("E2").Value = "$24.454.400.000"
$ws.Range("F2").Value = 44943.7729996875
$ws.Range("G2").Value = 44943.773196331

# Tweak column widths for "Persona a cargo" (B) and "Materiales" (C).
$ws.Columns.Item(2).ColumnWidth = 15.9
$ws.Columns.Item(3).ColumnWidth = 21.5
